$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as text (e.g. "30.145.65", "1.001").
# For updated prices that look like plain decimal numbers, force those
# specific cells to Text format first so Excel keeps them as strings
# instead of silently converting them to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.123.23'
$ws.Range('E2').Value = '  -4.55%  '
$ws.Range('D3').Value = '1.914.10'
$ws.Range('E3').Value = '  -3.92%  '
$ws.Range('D5').Value = '245.32'
$ws.Range('E5').Value = '  -3.48%  '
$ws.Range('D6').Value = '0.7009'
$ws.Range('E6').Value = '  -13.92%  '
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('D8').Value = '0.3218'
$ws.Range('E8').Value = '  -6.56%  '
$ws.Range('D9').Value = '25.86'
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').Value = '0.06843'
$ws.Range('E10').Value = '  -4.02%  '
$ws.Range('D11').Value = '0.7874'
$ws.Range('E11').Value = '  -7.17%  '
$ws.Range('D12').Value = '0.07929'
$ws.Range('E12').Value = '  -2.30%  '
$ws.Range('D13').Value = '1.918.26'
$ws.Range('E13').Value = '  -3.67%  '
$ws.Range('D14').Value = '5.357'
$ws.Range('E14').Value = '  -2.54%  '
$ws.Range('D15').Value = '93.45'
$ws.Range('E15').Value = '  -8.78%  '
$ws.Range('D16').Value = '259.60'
$ws.Range('E16').Value = '  -6.05%  '
$ws.Range('D17').Value = '14.32'
$ws.Range('E17').Value = '  +2.22%  '
$ws.Range('D18').Value = '30.142.82'
$ws.Range('E18').Value = '  -4.45%  '
$ws.Range('D19').Value = '5.784'
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('D20').Value = '0.000007834'
$ws.Range('E20').Value = '  -1.92%  '
$ws.Range('D21').Value = '2.169.78'
$ws.Range('E21').Value = '  -3.28%  '
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('D24').Value = '6.814'
$ws.Range('E24').Value = '  -1.55%  '
$ws.Range('D25').Value = '9.524'
$ws.Range('E25').Value = '  -1.94%  '
$ws.Range('D26').Value = '160.23'
$ws.Range('D27').Value = '18.70'
$ws.Range('E27').Value = '  -5.52%  '
$ws.Range('D28').Value = '0.1309'
$ws.Range('E28').Value = '  -16.03%  '
$ws.Range('D29').Value = '2.221'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '1.358'
$ws.Range('E30').Value = '  +0.50%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.548'
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('D32').Value = '4.390'
$ws.Range('E32').Value = '  -4.49%  '
$ws.Range('D33').Value = '4.162'
$ws.Range('E33').Value = '  -4.03%  '
$ws.Range('D34').Value = '0.05025'
$ws.Range('E34').Value = '  -3.42%  '
$ws.Range('D35').Value = '1.185'
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('D36').Value = '0.7396'
$ws.Range('E36').Value = '  -1.84%  '
$ws.Range('D37').Value = '2.724'
$ws.Range('E37').Value = '  -3.12%  '
$ws.Range('D38').Value = '0.01909'
$ws.Range('E38').Value = '  -4.95%  '
$ws.Range('D39').Value = '2.785'
$ws.Range('E39').Value = '  -5.15%  '
$ws.Range('D40').Value = '79.26'
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('D41').Value = '6.495'
$ws.Range('E41').Value = '  -2.54%  '
$ws.Range('E42').Value = '  -6.21%  '
$ws.Range('D43').Value = '2.003'
$ws.Range('E43').Value = '  -3.89%  '
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').Value = '0.8310'
$ws.Range('E45').Value = '  -3.13%  '
$ws.Range('D46').Value = '101.58'
$ws.Range('E46').Value = '  -4.79%  '
$ws.Range('D47').Value = '9.688'
$ws.Range('E47').Value = '  -2.91%  '
$ws.Range('D48').Value = '7.183'
$ws.Range('E48').Value = '  -4.48%  '
$ws.Range('D49').Value = '35.81'
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('D50').Value = '1.469'
$ws.Range('E50').Value = '  +1.81%  '
$ws.Range('D51').Value = '0.05910'
$ws.Range('E51').Value = '  -1.04%  '
